$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 486; this shifts rows 486:568 down to 487:569,
# which reproduces the historical (older) records moving down one slot as a
# new, more recent weekly record is inserted at the top of this block.
$ws.Rows(486).Insert()

# Populate the newly inserted row 486 with the new weekly record. Most of the
# descriptive/categorical fields repeat the constant values used throughout
# this block (market, region, category, etc.), mirroring the row that used to
# occupy this slot (now row 487) for K, L, M, N, O, P, Q, while D (date) and
# J (volume) hold the new observation's values.
$ws.Range("A486").Value = 10
$ws.Range("B486").Value = "Vega Modelo de Temuco"
$ws.Range("C486").Value = "La Araucanía"
$ws.Range("D486").Value = 45180
$ws.Range("E486").Value = 9
$ws.Range("F486").Value = 100112009
$ws.Range("G486").Value = "Acelga"
$ws.Range("H486").Value = "Sin especificar"
$ws.Range("I486").Value = "Primera"
$ws.Range("J486").Value = 90
$ws.Range("K486").Value = 8000
$ws.Range("L486").Value = 8000
$ws.Range("M486").Value = 8000
$ws.Range("N486").Value = "$/docena de atados (12 kilos)"
$ws.Range("O486").Value = "Provincia de Cautín"
$ws.Range("P486").Value = 667
$ws.Range("Q486").Value = 12
$ws.Range("R486").Value = "Hortaliza"

# Make sure the date cell keeps the date number format used by the rest of
# column D.
$ws.Range("D486").NumberFormat = $ws.Range("D487").NumberFormat
